# Update the dSF column (F) values for the miller_bryce worksheet.
# These values were repulled/recalculated; only specific rows change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -4
    5  = -2
    6  = 4
    8  = -1
    9  = 2
    10 = 1
    11 = 0
    12 = -2
    13 = 4
    14 = -1
    15 = 2
    16 = 4
    18 = 1
    19 = 6
    20 = -5
    21 = -4
    22 = 5
    23 = 1
    24 = 3
    25 = 3
    26 = -1
    27 = -3
    28 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
